$wb = $excel.ActiveWorkbook

# Row 58 on ALC (anchor G58=4606)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H58").Value = 5000
$ws.Range("J58").Value = 5000
$ws.Range("L58").Value = 15000
$ws.Range("N58").Value = -15300

# Row 94 on ALC (anchor G94=19905)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H94").Value = 2312.1428
$ws.Range("I94").Value = 2312.1428
$ws.Range("K94").Value = 2312.1428
$ws.Range("M94").Value = -1861.1428

# Row 98 on ALC (anchor G98=36237)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 1464.6538
$ws.Range("I98").Value = 1274.75
$ws.Range("K98").Value = 1274.75
$ws.Range("M98").Value = 223.25

# Row 107 on ALC (anchor G107=27766)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H107").Value = 1780.3334
$ws.Range("I107").Value = 846
$ws.Range("K107").Value = 846
$ws.Range("M107").Value = 1074

# Row 112 on ALC (anchor G112=27960)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 1425.8572
$ws.Range("J112").Value = 1441.6666
$ws.Range("L112").Value = 4324.9998
$ws.Range("N112").Value = -6540.9998

# Row 116 on ALC (anchor G116=27778)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 4149.8335
$ws.Range("I116").Value = 3299.6667
$ws.Range("K116").Value = 3299.6667
$ws.Range("M116").Value = 142.3332999999998

# Row 122 on ALC (anchor G122=36237)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 1464.6538
$ws.Range("I122").Value = 1274.75
$ws.Range("K122").Value = 3824.25
$ws.Range("M122").Value = -1374.25

# Row 132 on ALC (anchor G132=44049)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 1049.3793
$ws.Range("I132").Value = 940
$ws.Range("K132").Value = 2820
$ws.Range("M132").Value = -290

# Row 137 on ALC (anchor G137=44013)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1421
$ws.Range("I137").Value = 991.1667
$ws.Range("K137").Value = 2973.5001
$ws.Range("M137").Value = -423.5001000000002

# Row 141 on ALC (anchor G141=44161)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 4006759.2
$ws.Range("I141").Value = 5603418.5
$ws.Range("K141").Value = 16810255.5
$ws.Range("M141").Value = -16805075.5

# Row 45 on ARM (anchor G45=27714)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1730.6
$ws.Range("J45").Value = 1899.125
$ws.Range("L45").Value = 1899.125
$ws.Range("N45").Value = -2653.125

# Row 122 on ARM (anchor G122=36168)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 1250.2858
$ws.Range("I122").Value = 901.2857
$ws.Range("K122").Value = 2703.8571
$ws.Range("M122").Value = -253.8571000000002

# Row 132 on ARM (anchor G132=43997)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 1211.1628
$ws.Range("I132").Value = 1100.0256
$ws.Range("J132").Value = 2294.75
$ws.Range("K132").Value = 3300.0768
$ws.Range("L132").Value = 6884.25
$ws.Range("M132").Value = -770.0767999999998
$ws.Range("N132").Value = -11944.25

# Row 135 on ARM (anchor G135=42016)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H135").Value = 33099.715
$ws.Range("J135").Value = 33099.715
$ws.Range("L135").Value = 33099.715
$ws.Range("N135").Value = -43239.715

# Row 86 on BSM (anchor G86=12526)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 98218.664
$ws.Range("I86").Value = 3481.6
$ws.Range("K86").Value = 3481.6
$ws.Range("M86").Value = -2358.6

# Row 89 on BSM (anchor G89=12526)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 98218.664
$ws.Range("I89").Value = 3481.6
$ws.Range("K89").Value = 17408
$ws.Range("M89").Value = -11792

# Row 134 on BSM (anchor G134=43998)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2079.25
$ws.Range("I134").Value = 1786.2354
$ws.Range("J134").Value = 2790.8572
$ws.Range("K134").Value = 5358.706200000001
$ws.Range("L134").Value = 8372.571599999999
$ws.Range("M134").Value = -2823.706200000001
$ws.Range("N134").Value = -13442.5716

# Row 7 on CRP (anchor G7=5361)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 336.14285
$ws.Range("I7").Value = 467
$ws.Range("K7").Value = 467
$ws.Range("M7").Value = -354

# Row 16 on CRP (anchor G16=27691)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 770
$ws.Range("I16").Value = 849
$ws.Range("K16").Value = 849
$ws.Range("M16").Value = -562

# Row 22 on CRP (anchor G22=5367)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1190.7142
$ws.Range("I22").Value = 472.2
$ws.Range("J22").Value = 1589.8889
$ws.Range("K22").Value = 472.2
$ws.Range("L22").Value = 1589.8889
$ws.Range("M22").Value = -122.2
$ws.Range("N22").Value = -2289.8889

# Row 113 on CRP (anchor G113=27691)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 770
$ws.Range("I113").Value = 849
$ws.Range("K113").Value = 849
$ws.Range("M113").Value = 1321

# Row 132 on CRP (anchor G132=44019)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 1622.2727
$ws.Range("I132").Value = 1133.8
$ws.Range("J132").Value = 6507
$ws.Range("K132").Value = 3401.4
$ws.Range("L132").Value = 19521
$ws.Range("M132").Value = -871.3999999999996
$ws.Range("N132").Value = -24581

# Row 7 on CUL (anchor G7=4728)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 683.1667
$ws.Range("I7").Value = 300
$ws.Range("J7").Value = 718
$ws.Range("K7").Value = 900
$ws.Range("L7").Value = 2154
$ws.Range("M7").Value = -788
$ws.Range("N7").Value = -2378

# Row 63 on CUL (anchor G63=12866)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H63").Value = 0
$ws.Range("I63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("M63").ClearContents()

# Row 66 on CUL (anchor G66=12866)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H66").Value = 0
$ws.Range("I66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("M66").ClearContents()

# Row 86 on CUL (anchor G86=12892)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H86").Value = 350
$ws.Range("I86").Value = 200
$ws.Range("J86").Value = 400
$ws.Range("K86").Value = 600
$ws.Range("L86").Value = 1200
$ws.Range("M86").Value = 586
$ws.Range("N86").Value = -3572

# Row 89 on CUL (anchor G89=12892)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H89").Value = 350
$ws.Range("I89").Value = 200
$ws.Range("J89").Value = 400
$ws.Range("K89").Value = 1800
$ws.Range("L89").Value = 3600
$ws.Range("M89").Value = 4128
$ws.Range("N89").Value = -15456

# Row 132 on CUL (anchor G132=43972)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 2105
$ws.Range("I132").Value = 1225
$ws.Range("J132").Value = 2281
$ws.Range("K132").Value = 11025
$ws.Range("L132").Value = 20529
$ws.Range("M132").Value = -8495
$ws.Range("N132").Value = -25589

# Row 136 on CUL (anchor G136=44093)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H136").Value = 1104.5834
$ws.Range("I136").Value = 1104.5834
$ws.Range("K136").Value = 3313.7502
$ws.Range("M136").Value = 1786.2498

# Row 2 on GSM (anchor G2=5062)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 155
$ws.Range("J2").Value = 195
$ws.Range("L2").Value = 195
$ws.Range("N2").Value = -421

# Row 122 on GSM (anchor G122=36182)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2738.125
$ws.Range("I122").Value = 2881
$ws.Range("J122").Value = 2500
$ws.Range("K122").Value = 8643
$ws.Range("L122").Value = 7500
$ws.Range("M122").Value = -6193
$ws.Range("N122").Value = -12400

# Row 127 on GSM (anchor G127=34438)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H127").Value = 36666
$ws.Range("J127").Value = 36666
$ws.Range("L127").Value = 36666
$ws.Range("N127").Value = -46586

# Row 132 on GSM (anchor G132=44008)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3497758.2
$ws.Range("I132").Value = 3847284.2
$ws.Range("K132").Value = 11541852.6
$ws.Range("M132").Value = -11539322.6

# Row 7 on LTW (anchor G7=36249)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2562.2778
$ws.Range("I7").Value = 2448.2942
$ws.Range("K7").Value = 2448.2942
$ws.Range("M7").Value = -2336.2942

# Row 122 on LTW (anchor G122=36247)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 9041.944
$ws.Range("I122").Value = 9135
$ws.Range("J122").Value = 8800
$ws.Range("K122").Value = 27405
$ws.Range("L122").Value = 26400
$ws.Range("M122").Value = -24955
$ws.Range("N122").Value = -31300

# Row 126 on LTW (anchor G126=36249)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 2562.2778
$ws.Range("I126").Value = 2448.2942
$ws.Range("K126").Value = 7344.882599999999
$ws.Range("M126").Value = -4874.882599999999

# Row 132 on LTW (anchor G132=44058)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 1716.1111
$ws.Range("I132").Value = 1417.2727
$ws.Range("J132").Value = 2185.7144
$ws.Range("K132").Value = 4251.8181
$ws.Range("L132").Value = 6557.1432
$ws.Range("M132").Value = -1721.8181
$ws.Range("N132").Value = -11617.1432

# Row 107 on WVR (anchor G107=27746)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 702.7143
$ws.Range("I107").Value = 512.2727
$ws.Range("K107").Value = 1536.8181
$ws.Range("M107").Value = 383.1819

# Row 122 on WVR (anchor G122=36208)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 29453.178
$ws.Range("I122").Value = 57301.5
$ws.Range("J122").Value = 1604.8572
$ws.Range("K122").Value = 171904.5
$ws.Range("L122").Value = 4814.571599999999
$ws.Range("M122").Value = -169454.5
$ws.Range("N122").Value = -9714.571599999999

# Row 126 on WVR (anchor G126=36210)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 9761.076999999999
$ws.Range("I126").Value = 12090
$ws.Range("J126").Value = 1998
$ws.Range("K126").Value = 36270
$ws.Range("L126").Value = 5994
$ws.Range("M126").Value = -33800
$ws.Range("N126").Value = -10934

# Row 132 on WVR (anchor G132=44029)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1794.091
$ws.Range("J132").Value = 2606.2727
$ws.Range("L132").Value = 7818.8181
$ws.Range("N132").Value = -12878.8181
